# Update the dSF column (F) values on the active worksheet to reflect
# repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = -5
$ws.Range("F9").Value = -9
$ws.Range("F10").Value = -6
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = -5
$ws.Range("F19").Value = -2
